$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("D'Angelo Russell", "PG", "Los Angeles Lakers"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Draymond Green", "PF,C", "Golden State Warriors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
